$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("A7").Value2 = 112495160
$ws.Range("Q7").Value2 = 521847
$ws.Range("R7").Value2 = 7077543

# Row 8
$ws.Range("A8").Value2 = 112495159
$ws.Range("Q8").Value2 = 521837
$ws.Range("R8").Value2 = 7077579

# Row 12
$ws.Range("A12").Value2 = 112495164
$ws.Range("B12").Value2 = 56430
$ws.Range("E12").Value2 = 100109
$ws.Range("F12").Value2 = "Tretåig hackspett"
$ws.Range("G12").Value2 = "Picoides tridactylus"
$ws.Range("H12").Value2 = "(Linnaeus, 1758)"
$ws.Range("Q12").Value2 = 522103
$ws.Range("R12").Value2 = 7077453
$ws.Range("AC12").Value2 = "ringhack äldre"

# Row 13
$ws.Range("A13").Value2 = 112495168
$ws.Range("B13").Value2 = 56430
$ws.Range("E13").Value2 = 100109
$ws.Range("F13").Value2 = "Tretåig hackspett"
$ws.Range("G13").Value2 = "Picoides tridactylus"
$ws.Range("H13").Value2 = "(Linnaeus, 1758)"
$ws.Range("Q13").Value2 = 522175
$ws.Range("R13").Value2 = 7077257
$ws.Range("AC13").Value2 = "ringhack"

# Row 14
$ws.Range("A14").Value2 = 112495188
$ws.Range("B14").Value2 = 89549
$ws.Range("E14").Value2 = 1108
$ws.Range("F14").Value2 = "Harticka"
$ws.Range("G14").Value2 = "Pelloporus leporinus"
$ws.Range("H14").Value2 = "(Fr.) Krieglst."
$ws.Range("Q14").Value2 = 522460
$ws.Range("R14").Value2 = 7077294
$ws.Range("AC14").Value2 = ""

# Row 15
$ws.Range("A15").Value2 = 112495191
$ws.Range("B15").Value2 = 90799
$ws.Range("E15").Value2 = 1968
$ws.Range("F15").Value2 = "Grantaggsvamp"
$ws.Range("G15").Value2 = "Bankera violascens"
$ws.Range("H15").Value2 = "(Alb. & Schwein. : Fr.) Pouzar"
$ws.Range("Q15").Value2 = 522424
$ws.Range("R15").Value2 = 7077382
$ws.Range("AC15").Value2 = ""

# Row 16
$ws.Range("A16").Value2 = 112495170
$ws.Range("B16").Value2 = 56430
$ws.Range("E16").Value2 = 100109
$ws.Range("F16").Value2 = "Tretåig hackspett"
$ws.Range("G16").Value2 = "Picoides tridactylus"
$ws.Range("H16").Value2 = "(Linnaeus, 1758)"
$ws.Range("Q16").Value2 = 522164
$ws.Range("R16").Value2 = 7077256
$ws.Range("AC16").Value2 = "ringhack äldre"

# Row 17
$ws.Range("A17").Value2 = 112495166
$ws.Range("Q17").Value2 = 522179
$ws.Range("R17").Value2 = 7077220
$ws.Range("AC17").Value2 = "ringhack äldre"

# Row 18
$ws.Range("A18").Value2 = 112495185
$ws.Range("B18").Value2 = 89549
$ws.Range("E18").Value2 = 1108
$ws.Range("F18").Value2 = "Harticka"
$ws.Range("G18").Value2 = "Pelloporus leporinus"
$ws.Range("H18").Value2 = "(Fr.) Krieglst."
$ws.Range("Q18").Value2 = 522105
$ws.Range("R18").Value2 = 7077442
$ws.Range("AC18").Value2 = ""

# Row 19
$ws.Range("A19").Value2 = 112495165
$ws.Range("Q19").Value2 = 522086
$ws.Range("R19").Value2 = 7077322
$ws.Range("AC19").Value2 = "ringhack"

# Row 20
$ws.Range("A20").Value2 = 112495167
$ws.Range("Q20").Value2 = 522205
$ws.Range("R20").Value2 = 7077260

# Row 22
$ws.Range("A22").Value2 = 112495182
$ws.Range("B22").Value2 = 90826
$ws.Range("D22").Value2 = "LC"
$ws.Range("E22").Value2 = 4366
$ws.Range("F22").Value2 = "Skarp dropptaggsvamp"
$ws.Range("G22").Value2 = "Hydnellum peckii"
$ws.Range("H22").Value2 = "Banker"
$ws.Range("Q22").Value2 = 522424
$ws.Range("R22").Value2 = 7077374
$ws.Range("AC22").Value2 = ""

# Row 23
$ws.Range("A23").Value2 = 112495190
$ws.Range("B23").Value2 = 89571
$ws.Range("D23").Value2 = "NT"
$ws.Range("E23").Value2 = 5432
$ws.Range("F23").Value2 = "Granticka"
$ws.Range("G23").Value2 = "Porodaedalea chrysoloma"
$ws.Range("H23").Value2 = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q23").Value2 = 522469
$ws.Range("R23").Value2 = 7077316

# Row 24
$ws.Range("A24").Value2 = 112495173
$ws.Range("Q24").Value2 = 522099
$ws.Range("R24").Value2 = 7077313
$ws.Range("AC24").Value2 = "ringhack färska"

# Row 25
$ws.Range("A25").Value2 = 112495176
$ws.Range("B25").Value2 = 56430
$ws.Range("E25").Value2 = 100109
$ws.Range("F25").Value2 = "Tretåig hackspett"
$ws.Range("G25").Value2 = "Picoides tridactylus"
$ws.Range("H25").Value2 = "(Linnaeus, 1758)"
$ws.Range("Q25").Value2 = 522488
$ws.Range("R25").Value2 = 7077304
$ws.Range("AC25").Value2 = "ringhack äldre"

# Row 27
$ws.Range("A27").Value2 = 112495180
$ws.Range("B27").Value2 = 90857
$ws.Range("E27").Value2 = 5448
$ws.Range("F27").Value2 = "Svartvit taggsvamp"
$ws.Range("G27").Value2 = "Phellodon connatus"
$ws.Range("H27").Value2 = "(Schultz) nom.prov"
$ws.Range("Q27").Value2 = 522427
$ws.Range("R27").Value2 = 7077371
$ws.Range("AC27").Value2 = ""

# Row 29
$ws.Range("A29").Value2 = 112495169
$ws.Range("B29").Value2 = 56430
$ws.Range("E29").Value2 = 100109
$ws.Range("F29").Value2 = "Tretåig hackspett"
$ws.Range("G29").Value2 = "Picoides tridactylus"
$ws.Range("H29").Value2 = "(Linnaeus, 1758)"
$ws.Range("Q29").Value2 = 522163
$ws.Range("R29").Value2 = 7077259
$ws.Range("AC29").Value2 = "ringhack äldre"

# Row 30
$ws.Range("A30").Value2 = 112495172
$ws.Range("Q30").Value2 = 522109
$ws.Range("R30").Value2 = 7077299
$ws.Range("AC30").Value2 = "ringhack"

# Row 31
$ws.Range("A31").Value2 = 112495171
$ws.Range("Q31").Value2 = 522120
$ws.Range("R31").Value2 = 7077301
